# Update "Horarios Linea 141" workbook with the latest scrape (02:54:27).
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "LP1912" : refresh timestamps/minutes on existing rows and
# append the two new arrivals scraped in this run.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 02:54:27"
$ws1.Range("A3").Value = "Total filas: 5"

$ws1.Range("A6").Value = "02:54:27"
$ws1.Range("D6").Value = 8

$ws1.Range("A7").Value = "02:54:27"
$ws1.Range("D7").Value = 54

$ws1.Range("A8").Value = "02:54:27"
$ws1.Range("D8").Value = 67

$ws1.Range("A9").Value = "02:54:27"
$ws1.Range("B9").Value = "04:46"
$ws1.Range("C9").Value = "215_EL PELIGRO"
$ws1.Range("D9").Value = 112
$ws1.Range("E9").Value = "LP1912"

$ws1.Range("A10").Value = "02:54:27"
$ws1.Range("B10").Value = "04:53"
$ws1.Range("C10").Value = "11_ETCHEVERRY"
$ws1.Range("D10").Value = 119
$ws1.Range("E10").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "LP1912-215" : previously empty (Total filas: 0); now gains a
# header row plus the single matching arrival.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 02:54:27"
$ws2.Range("A3").Value = "Total filas: 1"

$ws2.Range("A5").Value = "Hora_Scrap"
$ws2.Range("B5").Value = "Hora_Llegada"
$ws2.Range("C5").Value = "Linea"
$ws2.Range("D5").Value = "Minutos"
$ws2.Range("E5").Value = "Parada"

$ws2.Range("A6").Value = "02:54:27"
$ws2.Range("B6").Value = "04:46"
$ws2.Range("C6").Value = "215_EL PELIGRO"
$ws2.Range("D6").Value = 112
$ws2.Range("E6").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "6203-6173" : still empty, only the refresh timestamp moves.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 02:54:27"
